# Generate Report for Archive
#
# Rows 8-10 on the Overview / zh-cn / de-de sheets describe three files:
#   014c6cdf-958b-4bc2-b78a-053e8c86b8b8
#   3f0a4fa7-07a4-48c0-9bc7-223a3d1bf54d
#   eda4a9c3-af38-4943-824e-7f400ff6e19a
#
# The refreshed status report re-sorts these three rows (3f0a4fa7, then
# eda4a9c3, then 014c6cdf) and flips the first two from
# "Ready for handoff" / "Ready for handoff" status to "In Translation",
# while the file that lands in row 10 keeps "Ready for handoff".

$wb = $excel.ActiveWorkbook

function Set-LinkDisplay($ws, $addr, $text) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $text
        }
    }
}

# ---------------------------------------------------------------------
# Overview sheet: columns A (file), B (zh-cn status), C (de-de status)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A8").Value = "3f0a4fa7-07a4-48c0-9bc7-223a3d1bf54d.md"
$wsOverview.Range("B8").Value = "In Translation"
$wsOverview.Range("C8").Value = "In Translation"
Set-LinkDisplay $wsOverview '$A$8' "3f0a4fa7-07a4-48c0-9bc7-223a3d1bf54d.md"

$wsOverview.Range("A9").Value = "eda4a9c3-af38-4943-824e-7f400ff6e19a.md"
$wsOverview.Range("B9").Value = "In Translation"
$wsOverview.Range("C9").Value = "In Translation"
Set-LinkDisplay $wsOverview '$A$9' "eda4a9c3-af38-4943-824e-7f400ff6e19a.md"

$wsOverview.Range("A10").Value = "014c6cdf-958b-4bc2-b78a-053e8c86b8b8.md"
$wsOverview.Range("B10").Value = "Ready for handoff"
$wsOverview.Range("C10").Value = "Ready for handoff"
Set-LinkDisplay $wsOverview '$A$10' "014c6cdf-958b-4bc2-b78a-053e8c86b8b8.md"

# ---------------------------------------------------------------------
# zh-cn sheet: columns A (file), C (status), D (latest handoff file)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A8").Value = "3f0a4fa7-07a4-48c0-9bc7-223a3d1bf54d.md"
$wsZhCn.Range("C8").Value = "In Translation"
$wsZhCn.Range("D8").Value = "3f0a4fa7-07a4-48c0-9bc7-223a3d1bf54d.d625f5e57c676e5a79e3c87bd9555d850b0e6cb6.zh-cn.xlf"
Set-LinkDisplay $wsZhCn '$A$8' "3f0a4fa7-07a4-48c0-9bc7-223a3d1bf54d.md"
Set-LinkDisplay $wsZhCn '$D$8' "3f0a4fa7-07a4-48c0-9bc7-223a3d1bf54d.d625f5e57c676e5a79e3c87bd9555d850b0e6cb6.zh-cn.xlf"

$wsZhCn.Range("A9").Value = "eda4a9c3-af38-4943-824e-7f400ff6e19a.md"
$wsZhCn.Range("C9").Value = "In Translation"
$wsZhCn.Range("D9").Value = "eda4a9c3-af38-4943-824e-7f400ff6e19a.630fd68283c0ce5ab6a872f896d824c96bc859de.zh-cn.xlf"
Set-LinkDisplay $wsZhCn '$A$9' "eda4a9c3-af38-4943-824e-7f400ff6e19a.md"
Set-LinkDisplay $wsZhCn '$D$9' "eda4a9c3-af38-4943-824e-7f400ff6e19a.630fd68283c0ce5ab6a872f896d824c96bc859de.zh-cn.xlf"

$wsZhCn.Range("A10").Value = "014c6cdf-958b-4bc2-b78a-053e8c86b8b8.md"
$wsZhCn.Range("C10").Value = "Ready for handoff"
$wsZhCn.Range("D10").Value = "014c6cdf-958b-4bc2-b78a-053e8c86b8b8.e44e71d4f0489edd6755148b97b69e11f7257c4a.zh-cn.xlf"
Set-LinkDisplay $wsZhCn '$A$10' "014c6cdf-958b-4bc2-b78a-053e8c86b8b8.md"
Set-LinkDisplay $wsZhCn '$D$10' "014c6cdf-958b-4bc2-b78a-053e8c86b8b8.e44e71d4f0489edd6755148b97b69e11f7257c4a.zh-cn.xlf"

# ---------------------------------------------------------------------
# de-de sheet: columns A (file), C (status), D (latest handoff file)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A8").Value = "3f0a4fa7-07a4-48c0-9bc7-223a3d1bf54d.md"
$wsDeDe.Range("C8").Value = "In Translation"
$wsDeDe.Range("D8").Value = "3f0a4fa7-07a4-48c0-9bc7-223a3d1bf54d.d625f5e57c676e5a79e3c87bd9555d850b0e6cb6.de-de.xlf"
Set-LinkDisplay $wsDeDe '$A$8' "3f0a4fa7-07a4-48c0-9bc7-223a3d1bf54d.md"
Set-LinkDisplay $wsDeDe '$D$8' "3f0a4fa7-07a4-48c0-9bc7-223a3d1bf54d.d625f5e57c676e5a79e3c87bd9555d850b0e6cb6.de-de.xlf"

$wsDeDe.Range("A9").Value = "eda4a9c3-af38-4943-824e-7f400ff6e19a.md"
$wsDeDe.Range("C9").Value = "In Translation"
$wsDeDe.Range("D9").Value = "eda4a9c3-af38-4943-824e-7f400ff6e19a.630fd68283c0ce5ab6a872f896d824c96bc859de.de-de.xlf"
Set-LinkDisplay $wsDeDe '$A$9' "eda4a9c3-af38-4943-824e-7f400ff6e19a.md"
Set-LinkDisplay $wsDeDe '$D$9' "eda4a9c3-af38-4943-824e-7f400ff6e19a.630fd68283c0ce5ab6a872f896d824c96bc859de.de-de.xlf"

$wsDeDe.Range("A10").Value = "014c6cdf-958b-4bc2-b78a-053e8c86b8b8.md"
$wsDeDe.Range("C10").Value = "Ready for handoff"
$wsDeDe.Range("D10").Value = "014c6cdf-958b-4bc2-b78a-053e8c86b8b8.e44e71d4f0489edd6755148b97b69e11f7257c4a.de-de.xlf"
Set-LinkDisplay $wsDeDe '$A$10' "014c6cdf-958b-4bc2-b78a-053e8c86b8b8.md"
Set-LinkDisplay $wsDeDe '$D$10' "014c6cdf-958b-4bc2-b78a-053e8c86b8b8.e44e71d4f0489edd6755148b97b69e11f7257c4a.de-de.xlf"
